$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing "getpadd" / "F" (row 5), shifting rows below it up.
$ws.Rows.Item(5).Delete()

# Update selection to match the post-edit state (row 5 selected as an entire row)
$ws.Range("A5:XFD5").Select()
